$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.04970888862596
$ws.Cells.Item(2, 4).Value = 1.05336706752593
$ws.Cells.Item(2, 5).Value = 1.04680113597715
$ws.Cells.Item(2, 6).Value = 1.061627691383946
$ws.Cells.Item(2, 9).Value = 1.038302735358962
$ws.Cells.Item(2, 10).Value = 1.054745901642773
$ws.Cells.Item(2, 11).Value = 1.056113232328338
$ws.Cells.Item(2, 12).Value = 1.049565555426539
$ws.Cells.Item(2, 13).Value = 1.064351237256879

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.051408124854707
$ws.Cells.Item(3, 4).Value = 1.054710776521811
$ws.Cells.Item(3, 5).Value = 1.048273350027508
$ws.Cells.Item(3, 6).Value = 1.063200686143314
$ws.Cells.Item(3, 9).Value = 1.038686550608598
$ws.Cells.Item(3, 10).Value = 1.056090861217224
$ws.Cells.Item(3, 11).Value = 1.057268612926101
$ws.Cells.Item(3, 12).Value = 1.050847771178331
$ws.Cells.Item(3, 13).Value = 1.065736980932223

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.052505325776799
$ws.Cells.Item(4, 4).Value = 1.055577957420256
$ws.Cells.Item(4, 5).Value = 1.049224047233271
$ws.Cells.Item(4, 6).Value = 1.064216531138524
$ws.Cells.Item(4, 9).Value = 1.038932505001579
$ws.Cells.Item(4, 10).Value = 1.056958514722632
$ws.Cells.Item(4, 11).Value = 1.058013399260951
$ws.Cells.Item(4, 12).Value = 1.051675023326692
$ws.Cells.Item(4, 13).Value = 1.066631172740828

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.052966046494881
$ws.Cells.Item(5, 4).Value = 1.055941981271521
$ws.Cells.Item(5, 5).Value = 1.049623269633532
$ws.Cells.Item(5, 6).Value = 1.064643126923598
$ws.Cells.Item(5, 9).Value = 1.039035332959586
$ws.Cells.Item(5, 10).Value = 1.057322657464131
$ws.Cells.Item(5, 11).Value = 1.058325840118521
$ws.Cells.Item(5, 12).Value = 1.052022228145743
$ws.Cells.Item(5, 13).Value = 1.067006507577721

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.053043372077187
$ws.Cells.Item(6, 4).Value = 1.056003071128438
$ws.Cells.Item(6, 5).Value = 1.049690274686088
$ws.Cells.Item(6, 6).Value = 1.064714727240329
$ws.Cells.Item(6, 9).Value = 1.039052564820106
$ws.Cells.Item(6, 10).Value = 1.057383762626949
$ws.Cells.Item(6, 11).Value = 1.058378261421791
$ws.Cells.Item(6, 12).Value = 1.052080492083909
$ws.Cells.Item(6, 13).Value = 1.067069494015104

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.05251148406419
$ws.Cells.Item(7, 4).Value = 1.055582823630303
$ws.Cells.Item(7, 5).Value = 1.049229383419868
$ws.Cells.Item(7, 6).Value = 1.064222233146748
$ws.Cells.Item(7, 9).Value = 1.038933881233967
$ws.Cells.Item(7, 10).Value = 1.056963382833779
$ws.Cells.Item(7, 11).Value = 1.058017576721558
$ws.Cells.Item(7, 12).Value = 1.051679664929409
$ws.Cells.Item(7, 13).Value = 1.066636190262356

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.050283641094669
$ws.Cells.Item(8, 4).Value = 1.053821659544516
$ws.Cells.Item(8, 5).Value = 1.047299081874204
$ws.Cells.Item(8, 6).Value = 1.062159710217081
$ws.Cells.Item(8, 9).Value = 1.038432946356401
$ws.Cells.Item(8, 10).Value = 1.055200986471818
$ws.Cells.Item(8, 11).Value = 1.056504287985325
$ws.Cells.Item(8, 12).Value = 1.04999939429882
$ws.Cells.Item(8, 13).Value = 1.064820074503731

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.046339576494525
$ws.Cells.Item(9, 4).Value = 1.050700333107754
$ws.Cells.Item(9, 5).Value = 1.043882467239554
$ws.Cells.Item(9, 6).Value = 1.058509587237848
$ws.Cells.Item(9, 9).Value = 1.03753170998789
$ws.Cells.Item(9, 10).Value = 1.052074875273023
$ws.Cells.Item(9, 11).Value = 1.053815690597578
$ws.Cells.Item(9, 12).Value = 1.047019543998641
$ws.Cells.Item(9, 13).Value = 1.061600442908638

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.043697091475469
$ws.Cells.Item(10, 4).Value = 1.048606828285333
$ws.Cells.Item(10, 5).Value = 1.041593897322076
$ws.Cells.Item(10, 6).Value = 1.056064959705235
$ws.Cells.Item(10, 9).Value = 1.036918224592824
$ws.Cells.Item(10, 10).Value = 1.049976374322297
$ws.Cells.Item(10, 11).Value = 1.052007980689353
$ws.Cells.Item(10, 12).Value = 1.045019627900092
$ws.Cells.Item(10, 13).Value = 1.059440367722662

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.042549580245487
$ws.Cells.Item(11, 4).Value = 1.047697197409005
$ws.Cells.Item(11, 5).Value = 1.040600213796262
$ws.Cells.Item(11, 6).Value = 1.055003601685753
$ws.Cells.Item(11, 9).Value = 1.036649529227562
$ws.Cells.Item(11, 10).Value = 1.049064140758654
$ws.Cells.Item(11, 11).Value = 1.051221472363606
$ws.Cells.Item(11, 12).Value = 1.044150346952669
$ws.Cells.Item(11, 13).Value = 1.058501656444982

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.042122832705173
$ws.Cells.Item(12, 4).Value = 1.047358838809663
$ws.Cells.Item(12, 5).Value = 1.040230695085916
$ws.Cells.Item(12, 6).Value = 1.054608929823063
$ws.Cells.Item(12, 9).Value = 1.036549261190091
$ws.Cells.Item(12, 10).Value = 1.04872474846828
$ws.Cells.Item(12, 11).Value = 1.050928753262771
$ws.Cells.Item(12, 12).Value = 1.043826949702147
$ws.Cells.Item(12, 13).Value = 1.058152456959303

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.042214394862573
$ws.Cells.Item(13, 4).Value = 1.047431439830911
$ws.Cells.Item(13, 5).Value = 1.04030997732361
$ws.Cells.Item(13, 6).Value = 1.054693608225575
$ws.Cells.Item(13, 9).Value = 1.03657079002456
$ws.Cells.Item(13, 10).Value = 1.048797574263854
$ws.Cells.Item(13, 11).Value = 1.050991568684906
$ws.Cells.Item(13, 12).Value = 1.043896342661316
$ws.Cells.Item(13, 13).Value = 1.058227385171448

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.042514315669563
$ws.Cells.Item(14, 4).Value = 1.047669238461959
$ws.Cells.Item(14, 5).Value = 1.040569677911028
$ws.Cells.Item(14, 6).Value = 1.054970986960327
$ws.Cells.Item(14, 9).Value = 1.036641250501196
$ws.Cells.Item(14, 10).Value = 1.049036097737302
$ws.Cells.Item(14, 11).Value = 1.051197287898914
$ws.Cells.Item(14, 12).Value = 1.044123625255199
$ws.Cells.Item(14, 13).Value = 1.05847280218991

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.042699038633472
$ws.Cells.Item(15, 4).Value = 1.047815689952704
$ws.Cells.Item(15, 5).Value = 1.040729631941973
$ws.Cells.Item(15, 6).Value = 1.055141830961384
$ws.Cells.Item(15, 9).Value = 1.036684602102939
$ws.Cells.Item(15, 10).Value = 1.049182986929536
$ws.Cells.Item(15, 11).Value = 1.051323961808072
$ws.Cells.Item(15, 12).Value = 1.044263593963278
$ws.Cells.Item(15, 13).Value = 1.058623942410754

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.043773177231833
$ws.Cells.Item(16, 4).Value = 1.048667130618234
$ws.Cells.Item(16, 5).Value = 1.041659786536608
$ws.Cells.Item(16, 6).Value = 1.056135338073073
$ws.Cells.Item(16, 9).Value = 1.036935992368214
$ws.Cells.Item(16, 10).Value = 1.050036840051349
$ws.Cells.Item(16, 11).Value = 1.052060098603827
$ws.Cells.Item(16, 12).Value = 1.045077248647073
$ws.Cells.Item(16, 13).Value = 1.059502594583707

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.044446062318364
$ws.Cells.Item(17, 4).Value = 1.049200371259041
$ws.Cells.Item(17, 5).Value = 1.042242511838029
$ws.Cells.Item(17, 6).Value = 1.05675777531013
$ws.Cells.Item(17, 9).Value = 1.037092862839226
$ws.Cells.Item(17, 10).Value = 1.050571475684805
$ws.Cells.Item(17, 11).Value = 1.05252084441102
$ws.Cells.Item(17, 12).Value = 1.045586740266394
$ws.Cells.Item(17, 13).Value = 1.060052835437792

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.044838227660146
$ws.Cells.Item(18, 4).Value = 1.049511100305223
$ws.Cells.Item(18, 5).Value = 1.042582144258262
$ws.Cells.Item(18, 6).Value = 1.057120561381946
$ws.Cells.Item(18, 9).Value = 1.037184068520771
$ws.Cells.Item(18, 10).Value = 1.050882976080775
$ws.Cells.Item(18, 11).Value = 1.052789227753086
$ws.Cells.Item(18, 12).Value = 1.045883600075099
$ws.Cells.Item(18, 13).Value = 1.060373455873156

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.04497189257511
$ws.Cells.Item(19, 4).Value = 1.049617000123167
$ws.Cells.Item(19, 5).Value = 1.042697906164856
$ws.Cells.Item(19, 6).Value = 1.057244216470794
$ws.Cells.Item(19, 9).Value = 1.037215117520278
$ws.Cells.Item(19, 10).Value = 1.050989131758015
$ws.Cells.Item(19, 11).Value = 1.052880678463466
$ws.Cells.Item(19, 12).Value = 1.04598476799142
$ws.Cells.Item(19, 13).Value = 1.060482724260236

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.044373901031187
$ws.Cells.Item(20, 4).Value = 1.049143190818798
$ws.Cells.Item(20, 5).Value = 1.042180018052264
$ws.Cells.Item(20, 6).Value = 1.056691021812783
$ws.Cells.Item(20, 9).Value = 1.037076062583217
$ws.Cells.Item(20, 10).Value = 1.050514149959735
$ws.Cells.Item(20, 11).Value = 1.05247144825293
$ws.Cells.Item(20, 12).Value = 1.045532109614007
$ws.Cells.Item(20, 13).Value = 1.059993833580549

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.042426010690994
$ws.Cells.Item(21, 4).Value = 1.04759922606441
$ws.Cells.Item(21, 5).Value = 1.040493214264008
$ws.Cells.Item(21, 6).Value = 1.054889317982487
$ws.Cells.Item(21, 9).Value = 1.036620514450325
$ws.Cells.Item(21, 10).Value = 1.048965873747048
$ws.Cells.Item(21, 11).Value = 1.051136724670758
$ws.Cells.Item(21, 12).Value = 1.044056710299003
$ws.Cells.Item(21, 13).Value = 1.058400547453964

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.041198331824135
$ws.Cells.Item(22, 4).Value = 1.046625684578137
$ws.Cells.Item(22, 5).Value = 1.03943021534984
$ws.Cells.Item(22, 6).Value = 1.053753984883721
$ws.Cells.Item(22, 9).Value = 1.036331414663115
$ws.Cells.Item(22, 10).Value = 1.047989233731396
$ws.Cells.Item(22, 11).Value = 1.050294199669673
$ws.Cells.Item(22, 12).Value = 1.043126125729924
$ws.Cells.Item(22, 13).Value = 1.057395768855433

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.041849433623998
$ws.Cells.Item(23, 4).Value = 1.047142045468666
$ws.Cells.Item(23, 5).Value = 1.039993966396084
$ws.Cells.Item(23, 6).Value = 1.054356090613272
$ws.Cells.Item(23, 9).Value = 1.036484927199071
$ws.Cells.Item(23, 10).Value = 1.048507274277905
$ws.Cells.Item(23, 11).Value = 1.050741157347508
$ws.Cells.Item(23, 12).Value = 1.043619728898149
$ws.Cells.Item(23, 13).Value = 1.057928710904593

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.044406508590583
$ws.Cells.Item(24, 4).Value = 1.049169029129155
$ws.Cells.Item(24, 5).Value = 1.04220825711172
$ws.Cells.Item(24, 6).Value = 1.056721185682366
$ws.Cells.Item(24, 9).Value = 1.037083654804856
$ws.Cells.Item(24, 10).Value = 1.05054005404878
$ws.Cells.Item(24, 11).Value = 1.052493769368765
$ws.Cells.Item(24, 12).Value = 1.045556795836962
$ws.Cells.Item(24, 13).Value = 1.060020494986667

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.047361463783018
$ws.Cells.Item(25, 4).Value = 1.051509453478215
$ws.Cells.Item(25, 5).Value = 1.044767606286391
$ws.Cells.Item(25, 6).Value = 1.059455158941106
$ws.Cells.Item(25, 9).Value = 1.037766917571452
$ws.Cells.Item(25, 10).Value = 1.052885547361144
$ws.Cells.Item(25, 11).Value = 1.054513417266077
$ws.Cells.Item(25, 12).Value = 1.047792217082301
$ws.Cells.Item(25, 13).Value = 1.062435156523321
